# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet right after "总计" (pushing every
# existing quarter tab one position to the right, which happens naturally
# because we insert rather than append), fills it with the new quarter's
# fund-holding figures, and updates the "总计" roll-up sheet with a new
# row for 2022-Q3 (the older rows keep their own name/value pairs, they
# simply end up one row lower).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating the existing
#    "2022-Q2" sheet (so header styling / column widths / formats match
#    the rest of the quarterly sheets) and drop it right after 总计.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($null, $totalSheet)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# New quarter's fund data (same fund 007280, new figures). These columns
# are stored as text (e.g. "1.12", not the number 1.12) everywhere else in
# the workbook, so force a text number-format before writing the values -
# otherwise Excel happily "helps" by coercing numeric-looking text back
# into a real number.
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.12"
$newSheet.Range("E2").Value = "90.06"
$newSheet.Range("F2").Value = "2.86"
$newSheet.Range("G2").Value = "0.0320"
$newSheet.Range("H2").Value = 5

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: add a row for 2022-Q3 and push the other
#    seven quarters down by one row. None of the existing quarter/value
#    pairs change - they simply land one row lower - and one brand new
#    row (2020-Q4) appears at the bottom, so just rewrite the whole
#    table (rows 2-9) with the new row included.
# ---------------------------------------------------------------------

# Row 9 is new, so clone the row-2 number formatting (bordered / bold /
# centered column-A style) onto it before filling in values.
$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A9:D9").PasteSpecial(-4122)

$quarters = @(
    @("2022-Q3", "0.03"),
    @("2022-Q2", "0.04"),
    @("2022-Q1", "0.05"),
    @("2021-Q4", "0.03"),
    @("2021-Q3", "0.03"),
    @("2021-Q2", "0.03"),
    @("2021-Q1", "0.05"),
    @("2020-Q4", "0.04")
)

for ($i = 0; $i -lt $quarters.Count; $i++) {
    $row = 2 + $i
    $label = $quarters[$i][0]
    $value = [double]$quarters[$i][1]

    $totalSheet.Cells.Item($row, 1).Value = $i
    $totalSheet.Cells.Item($row, 2).Value = $label
    $totalSheet.Cells.Item($row, 3).Value = 1
    $totalSheet.Cells.Item($row, 4).Value = $value
}
